$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 1008
$ws.Range("L3").Value = 1016
$ws.Range("L4").Value = 283
$ws.Range("L5").Value = 67
$ws.Range("K6").Value = 9121
$ws.Range("L6").Value = 1047
$ws.Range("K7").Value = 27533
$ws.Range("L7").Value = 3421

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 58
$ws.Range("L3").Value = 65
$ws.Range("L7").Value = 207

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L3").Value = 35
$ws.Range("L4").Value = 2
$ws.Range("L7").Value = 75

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 31
$ws.Range("L3").Value = 56
$ws.Range("L7").Value = 149

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L3").Value = 29
$ws.Range("L7").Value = 118

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 24
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L4").Value = 6
$ws.Range("L7").Value = 54

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L7").Value = 110
$ws.Range("L8").Value = 207
$ws.Range("L10").Value = 23
$ws.Range("L11").Value = 52
$ws.Range("L13").Value = 3
$ws.Range("L15").Value = 23
$ws.Range("K20").Value = 673
$ws.Range("L20").Value = 88
$ws.Range("L22").Value = 10
$ws.Range("L30").Value = 20
$ws.Range("L31").Value = 37
$ws.Range("L33").Value = 149
$ws.Range("L37").Value = 118
$ws.Range("L42").Value = 109
$ws.Range("L48").Value = 54
$ws.Range("L51").Value = 46
$ws.Range("L52").Value = 64
$ws.Range("L55").Value = 38
$ws.Range("L57").Value = 18
$ws.Range("L60").Value = 20
$ws.Range("L63").Value = 16
$ws.Range("L64").Value = 28
$ws.Range("L65").Value = 70
$ws.Range("L66").Value = 3
$ws.Range("L67").Value = 126
$ws.Range("L77").Value = 20
$ws.Range("L78").Value = 53
$ws.Range("L79").Value = 100
$ws.Range("L80").Value = 13
$ws.Range("L82").Value = 11
$ws.Range("L83").Value = 75
$ws.Range("L84").Value = 33
$ws.Range("L85").Value = 175
$ws.Range("L86").Value = 25
$ws.Range("L87").Value = 12
$ws.Range("L90").Value = 29
$ws.Range("L91").Value = 44
$ws.Range("L92").Value = 9
$ws.Range("L93").Value = 19
$ws.Range("L96").Value = 30
$ws.Range("L98").Value = 26
$ws.Range("L99").Value = 54
$ws.Range("K101").Value = 27533
$ws.Range("L101").Value = 3421

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L2").Value = 11
$ws.Range("L7").Value = 37

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 40
$ws.Range("L3").Value = 32
$ws.Range("L7").Value = 126

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L2").Value = 13
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 33

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 54

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 29
$ws.Range("L6").Value = 47
$ws.Range("L7").Value = 109

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("L5").Value = 2
$ws.Range("L6").Value = 3

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 23

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L3").Value = 12
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 53

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L2").Value = 17
$ws.Range("L7").Value = 38

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L2").Value = 16
$ws.Range("L7").Value = 30

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 44

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 34
$ws.Range("L5").Value = 5
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 100

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L5").Value = 2
$ws.Range("L7").Value = 28

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 26
$ws.Range("K6").Value = 195
$ws.Range("K7").Value = 673
$ws.Range("L7").Value = 88

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L4").Value = 1
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 19

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 110

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 8
$ws.Range("L7").Value = 23

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L4").Value = 1
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L3").Value = 1
$ws.Range("L7").Value = 3

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 17
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 52

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("L2").Value = 4
$ws.Range("L7").Value = 9

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L2").Value = 4
$ws.Range("L7").Value = 25

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 15
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L2").Value = 15
$ws.Range("L7").Value = 46

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("L4").Value = 2
$ws.Range("L7").Value = 18

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L2").Value = 8
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 76
$ws.Range("L4").Value = 14
$ws.Range("L7").Value = 175

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L2").Value = 5
$ws.Range("L7").Value = 10

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("L5").Value = 5
$ws.Range("L6").Value = 11

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L3").Value = 9
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 13

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 22
$ws.Range("L4").Value = 8
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 64

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("L2").Value = 2
$ws.Range("L7").Value = 12
